$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '66.239.97'
$ws.Range("E2").Value = '  -0.58%  '

$ws.Range("D3").Value = '3.209.89'
$ws.Range("E3").Value = '  +0.53%  '

$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").Value = '607.59'
$ws.Range("E5").Value = '  +1.77%  '

$ws.Range("D6").Value = '156.70'
$ws.Range("E6").Value = '  +1.14%  '

$ws.Range("E7").Value = '  +0.09%  '

$ws.Range("D8").Value = '3.210.55'
$ws.Range("E8").Value = '  +0.93%  '

$ws.Range("D9").Value = '0.553'
$ws.Range("E9").Value = '  -0.41%  '

$ws.Range("D10").Value = '0.160'
$ws.Range("E10").Value = '  -0.10%  '

$ws.Range("D11").Value = '5.67'
$ws.Range("E11").Value = '  -4.77%  '

$ws.Range("E12").Value = '  -2.66%  '

$ws.Range("D13").Value = '0.0000270'
$ws.Range("E13").Value = '  +0.30%  '

$ws.Range("D14").Value = '38.54'
$ws.Range("E14").Value = '  -1.97%  '

$ws.Range("D15").Value = '3.737.39'
$ws.Range("E15").Value = '  +0.65%  '

$ws.Range("D16").Value = '66.432.39'
$ws.Range("E16").Value = '  -0.08%  '

$ws.Range("D17").Value = '7.37'
$ws.Range("E17").Value = '  -1.58%  '

$ws.Range("D18").Value = '3.209.23'
$ws.Range("E18").Value = '  +0.58%  '

$ws.Range("E19").Value = '  +1.42%  '

$ws.Range("D20").Value = '510.26'
$ws.Range("E20").Value = '  -1.44%  '

$ws.Range("D21").Value = '15.29'
$ws.Range("E21").Value = '  -0.84%  '

$ws.Range("D22").Value = '0.732'
$ws.Range("E22").Value = '  -0.94%  '

$ws.Range("D23").Value = '8.02'
$ws.Range("E23").Value = '  -0.80%  '

$ws.Range("D24").Value = '14.63'
$ws.Range("E24").Value = '  -2.30%  '

$ws.Range("D25").Value = '85.19'
$ws.Range("E25").Value = '  -0.79%  '

$ws.Range("E26").Value = '  -0.13%  '

$ws.Range("D27").Value = '3.00'
$ws.Range("E27").Value = '  -0.27%  '

$ws.Range("D28").Value = '9.05'
$ws.Range("E28").Value = '  -2.76%  '

$ws.Range("D29").Value = '2.36'
$ws.Range("E29").Value = '  +0.16%  '

$ws.Range("D30").Value = '0.129'
$ws.Range("E30").Value = '  +42.65%  '

$ws.Range("D31").Value = '2.93'
$ws.Range("E31").Value = '  -0.51%  '

$ws.Range("D32").Value = '6.98'
$ws.Range("E32").Value = '  -2.02%  '

$ws.Range("D33").Value = '28.19'
$ws.Range("E33").Value = '  -0.63%  '

$ws.Range("E34").Value = '  +0.40%  '

$ws.Range("E35").Value = '  -4.55%  '

$ws.Range("D36").Value = '6.51'
$ws.Range("E36").Value = '  -0.66%  '

$ws.Range("D37").Value = '500.78'
$ws.Range("E37").Value = '  -2.50%  '

$ws.Range("D38").Value = '55.35'
$ws.Range("E38").Value = '  +0.71%  '

$ws.Range("D39").Value = '0.0₃0771'
$ws.Range("E39").Value = '  +13.75%  '

$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = '0.0422'
$ws.Range("E40").Value = '  -1.02%  '

$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").Value = '0.131'
$ws.Range("E41").Value = '  +1.40%  '

$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").Value = '3.03'
$ws.Range("E42").Value = '  +4.18%  '

$ws.Range("D43").Value = '8.74'
$ws.Range("E43").Value = '  -2.07%  '

$ws.Range("D44").Value = '0.298'
$ws.Range("E44").Value = '  -1.89%  '

$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '2.921.39'
$ws.Range("E45").Value = '  +0.35%  '

$ws.Range("B46").Value = 'Fetch.AI'
$ws.Range("C46").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D46").Value = '2.44'
$ws.Range("E46").Value = '  -0.53%  '

$ws.Range("D47").Value = '28.16'
$ws.Range("E47").Value = '  -2.10%  '

$ws.Range("E48").Value = '  +2.80%  '

$ws.Range("E49").Value = '  -0.05%  '

$ws.Range("E50").Value = '  -0.73%  '

$ws.Range("D51").Value = '122.02'
$ws.Range("E51").Value = '  -0.44%  '
